$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.15
$ws.Range("C2").Value = 0.9399999999999999
$ws.Range("D2").Value = 69.66

$ws.Range("B3").Value = 0.88
$ws.Range("C3").Value = -0.78
$ws.Range("D3").Value = 35.03

$ws.Range("B4").Value = 1.03
$ws.Range("C4").Value = 0.35
$ws.Range("D4").Value = 61.08

$ws.Range("B5").Value = 0.6899999999999999
$ws.Range("C5").Value = -1.33
$ws.Range("D5").Value = 39.46

$ws.Range("B6").Value = 0.9
$ws.Range("C6").Value = -0.54
$ws.Range("D6").Value = 48.47

$ws.Range("B7").Value = 0.91
$ws.Range("C7").Value = -0.68
$ws.Range("D7").Value = 38.04
